$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 15 play-by-play yardage data ---
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = '1 11 5 15 3 3 6 2 8 2 10 5 3 9 5 9 3 11 1 1 6 6 14 2 -1 3 2 1 1 1 2 7 4 0 2 6 3 2 2 0 0 4 2 12 0 1 7 1 5 9 -7 8 6 1 0 5 0 7 5 3 6 5 3 3 3 3 8 5 2 7 9 5 7 5 -1 1 -2 1 -3 3 11 3 14 6 3 2 5 5 4 6 10 6 3 3 3 2 5 4 2 3 8 1 3 2 5 8 1 8 5 3 4 6 9 2 9 4 4 4 8 7 1 0 3 2 8 1 5 2 1 5 3 1 23 2 10 5 3 12 3 3 0 11 2 1 8 3 -2 1 5 4 2 3 1 3 4 3 2 0 0 2 18 3 -1 6 2 5 -3 3 9 1 7 1 12 4 -2 20 10 1 4 3 13 3 1 3 -1 3 3 11 7 1 1 3 2 15 2 4 2 3 6 5 1 22 1 3 9 2 9 11 9 9 14 2 4 1 3 8 -1 3 1 6 -1 3 2 44 2 4 2 2 3 1 2 5 10 3 1 6 8 0 2 2 6 4 4 6 13 3 0 5 3 2 6 -1 4 2 14 9 1 -2 4 3 1 9 1 -1 28 3 1 17 7 3 1 2 -2 3 4 8 2 0 0 2 2 1 4 -5 4 2 5 3 -1 3 2 9 9 3 -1 25 24 7 0 4 8 6 5 9 23 4 6 8 0 7 7 5 11 7 1 1 3 10 -1 2 4 7 1 3 6 0 8 9 2 16 5 1 0 -2 7 13 -2 2 2 -1 9 3 0 9 3 8 3 1 -1 2 -1 0 0 44'   # OFF R
$ydsWs.Range("B3").Value = '14 3 17 2 3 3 8 5 6 4 10 8 55 9 7 4 16 3 7 4 27 3 8 23 9 12 17 11 11 6 13 4 9 7 10 12 7 -3 11 9 6 18 2 11 6 5 6 19 15 58 9 21 15 3 16 8 4 11 72 20 11 10 12 49 14 -4 11 12 6 32 19 13 8 6 23 7 14 5 9 14 6 31 15 -5 6 7 13 8 11 13 6 4 5 18 7 16 5 14 12 14 -1 15 1 38 2 14 5 9 4 16 4 -1 9 2 19 5 14 7 19 11 20 4 15 -1 14 19 9 7 6 5 9 25 5 7 22 8 8 10 8 6 7 6 46 11 3 12 21 21 8 10 34 26 22 5 17 15 5 7 6 18 19 5 10 9 26 21 -6 17 2 5 2 10 1 20 6 5 22 4 7 8 14 13 9 7 10 17 9 24 11 -3 7 6 3 7 6 12 -2 -6 25 27 9 24 9 1 9 30 6 1 5 13 13 70 19 4 2 15 13 6 8 6 18 19 6 21 7 26 5'   # OFF P
$ydsWs.Range("C2").Value = '-1 1 6 6 3 3 5 3 3 0 3 5 3 2 6 1 1 0 2 5 -1 -1 16 6 0 7 3 5 9 0 2 -2 0 11 0 -1 1 3 2 -5 3 7 -2 4 3 6 1 2 0 1 17 17 2 2 2 0 -1 7 12 0 4 0 1 7 8 1 17 20 1 1 -3 2 4 1 4 2 1 6 6 5 8 3 0 1 8 5 1 -7 13 5 3 3 1 4 1 2 3 3 2 5 6 16 1 9 2 11 2 -2 -2 16 2 1 3 6 0 0 1 8 12 3 -3 1 1 5 12 2 14 3 6 1 4 5 0 1 -1 5 4 -1 -1 -1 4 13 9 -1 2 -1 0 24 6 2 5 6 3 1 0 4 1 5 0 2 2 1 3 4 4 5 1 2 1 0 1 3 -4 -1 2 1 2 4 2 5 4 -1 6 5 7 1 2 1 5 3 2 2 1 1 0 2 3 4 5 6 -5 0 1 2 3 4 -3 4 1 25 4 4 20 4 3 3 2 4 4 18 9 4 13 4 4 4 5 6 0 2 5 6 0 4 -3 7 5 10 6 2 -6 3 -1 7 2 1 14 1 3 24 6 6 9 5 2 8 -1 5 3 3 -3 1 2 4 2 3 -1 4 10 1 11 15 0 -5 9 6 3 2 -2 0 3 2 1 2 33 2 1 3 4 -2 4 58 2 10 2 3 5 3 3 10 3 2 18 3 1 2 -1 7 2 1 1 14 3 7 2 8 5'   # DEF R
$ydsWs.Range("C3").Value = '4 6 7 14 4 31 10 19 7 8 5 6 2 7 32 20 20 19 8 5 13 18 19 5 16 7 17 2 10 6 7 13 8 1 4 8 11 11 16 9 -1 14 7 2 9 6 17 31 13 7 3 10 8 27 7 9 4 4 4 4 11 5 22 9 10 4 5 -5 10 13 -4 6 18 4 27 51 19 7 8 17 20 8 4 5 21 1 54 8 5 11 3 28 17 8 18 1 23 20 9 18 8 7 17 13 7 16 10 7 -7 9 3 6 4 6 9 26 27 84 5 -4 5 10 10 1 7 12 9 28 0 8 25 12 7 8 8 4 20 31 4 9 7 44 6 7 16 6 6 41 9 7 -1 8 2 12 8 50 11 34 7 39 13 6 0 8 7 5 16 7 20 7 6 3 12 49 16 13 -4 9 64 3 2 50 9 16 17 10 1 16 16 9 10 -1 2 5 7 2 39 0 8 4 16 14 33 -2 4 14 6 14 14 -1 23 14 9 7 7 26 1 6 12 7 9 4 17 4 7 14 5 21 4 5 6 28 23 9 5 1 41 7 1 21 10 10 5 17 1 4 20 0 12 2 10 1 3 4 13 25 9 4 3 8 3 2 22 16 11 9 1 6 6 27 3 17 14 14 8 13'   # DEF P

# --- OFF sheet: updated season totals ---
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("C2").Value = 165
$offWs.Range("D2").Value = 4
$offWs.Range("E2").Value = 10
$offWs.Range("F2").Value = 57
$offWs.Range("G2").Value = 60
$offWs.Range("H2").Value = 7
$offWs.Range("I2").Value = 11
$offWs.Range("J2").Value = 25
$offWs.Range("N2").Value = 16
$offWs.Range("O2").Value = 17
$offWs.Range("P2").Value = 10
$offWs.Range("B3").Value = 11
$offWs.Range("C3").Value = 133
$offWs.Range("E3").Value = 32
$offWs.Range("F3").Value = 79
$offWs.Range("H3").Value = 24
$offWs.Range("I3").Value = 44
$offWs.Range("J3").Value = 36
$offWs.Range("L3").Value = 203
$offWs.Range("M3").Value = 114
$offWs.Range("Q3").Value = 408

# --- DEF sheet: updated season totals ---
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 157
$defWs.Range("F2").Value = 48
$defWs.Range("G2").Value = 37
$defWs.Range("H2").Value = 5
$defWs.Range("N2").Value = 21
$defWs.Range("O2").Value = 16
$defWs.Range("P2").Value = 6
$defWs.Range("B3").Value = 10
$defWs.Range("C3").Value = 135
$defWs.Range("E3").Value = 29
$defWs.Range("F3").Value = 92
$defWs.Range("G3").Value = 21
$defWs.Range("H3").Value = 23
$defWs.Range("I3").Value = 59
$defWs.Range("J3").Value = 45
$defWs.Range("L3").Value = 250
$defWs.Range("M3").Value = 140
$defWs.Range("Q3").Value = 478

# --- ST sheet: updated totals and appended Week 15 per-game data ---
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 62
$stWs.Range("D2").Value = 59
$stWs.Range("F2").Value = 80
$stWs.Range("G2").Value = 77
$stWs.Range("J2").Value = 40
$stWs.Range("K2").Value = 36
$stWs.Range("B3").Value = 39
$stWs.Range("B4").Value = '57 65 64 66 56 65 60 63 65 63 65 61 62 47 64 63 59 68 69 65 68 59 62'
$stWs.Range("B5").Value = '17 15 22 26 17 23 13 17 19 21 25 27 20 7 18 14 20 22 30 23 21 23 42'
$stWs.Range("D3").Value = '45 58 44 60 42 52 59 37 51 44 54 56 44 34 55 43 46 60 57 52 53 58 41 41 48 41 56 45 50 63 39 45 38 39 38 41 41 43 62 53 39 54 49 42 52 42 43 53 47 59 52 55 33 62 45 62 48 49 54'
$stWs.Range("D4").Value = '0 17 0 10 14 0 0 0 0 8 0 18 12 0 26 0 0 0 0 10 0 8 0 0 0 0 11 0 3 11 0 0 0 0 0 9 0 0 12 10 0 8 0 13 11 10 1 18 7 16 9 0 0 12 0 28 5 0 11'
$stWs.Range("D5").Value = '0 0 13 0 0 0 0 0 25 0 0 0 0 9 0 0 0 13 14 2 0 2 0 0 0 0 0 12 0 0 17 7 26 0 4 0 0 0 8 0 14 7 12 0 8 19 0 8 0 2 0 11 12 13 0 0 18'

# --- TURNS sheet: updated Road fumbles lost ---
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("D3").Value = 4

# --- PEN sheet: updated penalty counts ---
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("D3").Value = 3
$penWs.Range("D4").Value = 9

Write-Host "Week 15 simulations added"
